$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A1").Value = 11046
$ws.Range("B1").Value = "picanha"
$ws.Range("C1").Value = 800

$ws.Range("A2:C11").ClearContents()
$ws.Rows.Item(12).Delete()

[void]$ws.Range("C1").Select()
